$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 27: Results changes from PASS to SKIP -----------------------------
$ws.Range("E27").Value2 = "SKIP"

# --- Row 28: new test case, formatted like row 27 --------------------------
$ws.Range("A27:E27").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)

$ws.Range("A28").Value2 = "ProfileTypeaheadCountySelectTest"
$ws.Range("B28").Value2 = "TBD"
$ws.Range("C28").Value2 = "Verify that user is able to add 'country' using typeahead"
$ws.Range("D28").Value2 = "Y"
$ws.Range("E28").Value2 = "PASS"

$italicPart = $ws.Range("C28").Characters(34, 8)
$italicPart.Font.Italic = $true

$tailPart = $ws.Range("C28").Characters(42, 17)
$tailPart.Font.Size = 11

# --- Selection ---------------------------------------------------------------
$ws.Range("C18").Select()
